# Update API Runtime Excel files:
#  - Rename the data sheet "Sheet1" -> "Data" and make it the active sheet/tab
#  - Turn the Legend sheet into a small table (Table1) with a proper header
#    row ("Column1"/"Column2") above the existing legend key/value rows
#  - Size the legend columns and refresh the selection/active tab bookkeeping

$wb = $excel.ActiveWorkbook

# --- Sheet1 -> Data -------------------------------------------------------
$wsData = $wb.Worksheets.Item(1)
$wsData.Name = "Data"

# --- Legend sheet: insert header row, add table ---------------------------
$wsLegend = $wb.Worksheets.Item(2)

# Push the existing 6 rows of legend data down one row and add headers
$wsLegend.Rows.Item(1).Insert()
$wsLegend.Range("A1").Value = "Column1"
$wsLegend.Range("B1").Value = "Column2"

# Match the column widths used for the new table
$wsLegend.Range("A1:B7").ColumnWidth = 10.43

# Turn A1:B7 into an actual Excel Table (ListObject)
$lo = $wsLegend.ListObjects.Add(1, $wsLegend.Range("A1:B7"), $null, 1)
$lo.Name = "Table1"

# Select the whole table range on the Legend sheet (no longer the active tab)
[void]$wsLegend.Range("A1:B7").Select()

# --- Make "Data" the active/selected tab again -----------------------------
$wsData.Activate()
